# Adds a new article row (Savora mustard-based dressing) to the "Maestro"
# articles sheet. In the source data this new record lands at row 11,
# pushing every existing row from 11 down through 42 down by one (to
# 12-43). Rows are shifted by copying whole rows bottom-up (so we never
# clobber data we still need) which preserves each row's existing cell
# formatting/styles exactly, then the freshly freed row 11 is populated
# with the new article's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 42; $r -ge 11; $r--) {
    $ws.Rows.Item($r).Copy()
    $ws.Rows.Item($r + 1).PasteSpecial(-4104)
}

# Row 43 didn't exist before this edit, so the paste above landed on a
# brand-new row; make sure column A keeps the integer display format
# ("0") used by every Codigo cell in this sheet.
$ws.Cells.Item(43, 1).NumberFormat = "0"

# Column O carries a special "ImagenExactaDelArticulo" style (bordered,
# right-aligned, wrapped text) for rows 2-26 only; from row 27 on the
# cells use the plain default style. The whole-row PasteSpecial above
# doesn't reliably re-apply that style when the destination row
# previously had no explicit style of its own (exactly the old
# row25->row26 boundary, now shifted down to row26->row27), so re-copy
# that one cell's formatting explicitly - row 25 already carries the
# correct style at this point (it received old row 24's formatted data).
$ws.Cells.Item(25, 15).Copy()
$ws.Cells.Item(26, 15).PasteSpecial(-4122)

$ws.Cells.Item(11, 1).Value = 7794000006478
$ws.Cells.Item(11, 2).Value = "Aderezo a base de"
$ws.Cells.Item(11, 3).Value = "mostaza"
$ws.Cells.Item(11, 4).Value = "original"
$ws.Cells.Item(11, 5).Value = "Savora"
$ws.Cells.Item(11, 6).Value = 250
$ws.Cells.Item(11, 7).Value = "gr."
$ws.Cells.Item(11, 8).Value = "Pouch"
$ws.Cells.Item(11, 9).Value = "Aderezos"
$ws.Cells.Item(11, 10).Value = "Argentina"
$ws.Cells.Item(11, 11).Value = 6
$ws.Cells.Item(11, 12).Value = $false
$ws.Cells.Item(11, 13).Value = $true
$ws.Cells.Item(11, 14).Value = "C:\VentaSoft\Imágenes de artículos\7794000006478.png"
$ws.Cells.Item(11, 15).Value = $false
